# Rename the "_old" / "_new" suffixed column headers in row 1 to the
# respective input-file-name based suffixes ("_FV2310" / "_FV2404"),
# then turn the header/data range into a proper Excel Table and freeze
# the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
# Columns A-J (1-10) carry the "_old" suffix -> "_FV2310"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace '_old$', '_FV2310')
}

# Column K (11) is "diff" - untouched

# Columns L-U (12-21) carry the "_new" suffix -> "_FV2404"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value() -replace '_new$', '_FV2404')
}

# --- 2. Turn A1:U70 into an Excel Table (Table1) ----------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
